$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "SVM" row label to "SVM_local"
$ws.Range("A2").Value = "SVM_local"

# Bring over the row-2 formatting (bold/border/centered style + the
# trailing blank "ks" cell) onto the new row 3 before filling it in, so
# the new row reuses the same cell style as row 2 instead of creating a
# new one.
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)

# Populate the new AdaBoostClassifier_local row of results
$ws.Range("A3").Value = "AdaBoostClassifier_local"
$ws.Range("B3").Value = 84.02777777777779
$ws.Range("C3").Value = 91.66666666666666
$ws.Range("D3").Value = 91.66666666666666
$ws.Range("E3").Value = 87.68115942028986
$ws.Range("F3").Value = 0.6498316498316499
